$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Before insert D7:" $ws.Range("D7").Value2
$ws.Columns("D").Insert()
Write-Host "After insert D7:" $ws.Range("D7").Value2
Write-Host "After insert E7:" $ws.Range("E7").Value2
